# Updated templates for Victim Attorney in COS.
#
# 1) Turn the section's single default header/footer into the full
#    even / default / first trio (mirrors Word's "Different Odd & Even
#    Pages" + "Different First Page" split): unlink the even header from
#    the previous section and touch its Range so Word materialises
#    separate header1/2/3.xml + footer1/2/3.xml parts, with the original
#    content preserved on the "default" (primary) header/footer.
# 2) Append a new "Victim's Attorney" line to the default footer's
#    "Prosecutor's Office ..." paragraph.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

$hEven = $sec.Headers(3)
$hEven.LinkToPrevious = $false
$hEven.Range.Text = ""

$footerDefault = $sec.Footers(1)

$searchRange = $footerDefault.Range.Duplicate
$marker = "County Jail: PS   EM;"
$found = $searchRange.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $endPos = $searchRange.End
    $apostrophe = [char]0x2019

    $ip1 = $footerDefault.Range.Duplicate
    $ip1.SetRange($endPos, $endPos)
    $ip1.InsertAfter(" ")

    $ip2 = $footerDefault.Range.Duplicate
    $ip2.SetRange($endPos + 1, $endPos + 1)
    $ip2.InsertAfter("Victim" + $apostrophe + "s Attorney (if applicable): PS   OS   EM")
}

Write-Output "done"
